$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated s_vals data (regenerated to filter save games)
# Columns: B=TB, C=d2S, D=K, E=IP, G=sum (F=Win unchanged, A=date unchanged)

$data = @{
    2  = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 6.15379541431027 }
    3  = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.1496068669990043; E = 0.5333859586016987;  G = 5.582307763322248 }
    4  = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.1496068669990043; E = 0.5333859586016987;  G = 5.582307763322248 }
    5  = @{ B = 1.445647641019636;  C = 0.3048912486333797; D = 0.1496068669990043; E = 0.5333859586016987;  G = 2.433531715253719 }
    6  = @{ B = 0.6545652718822623; C = 0.3048912486333797; D = 0.1496068669990043; E = 0.5333859586016987;  G = 1.642449346116345 }
    7  = @{ B = 1.445647641019636;  C = 1.626987699542094;  D = 0.7210945179870265; E = 13.86384647080068;   G = 17.65757632934944 }
    8  = @{ B = 0.1169995834814548; C = 0.04103571897497393;D = 0.7210945179870265; E = 0.5333859586016987;  G = 1.412515779045154 }
    9  = @{ B = 0.6545652718822623; C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 3.536033448013082 }
    10 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 6.15379541431027 }
    11 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 6.15379541431027 }
    12 = @{ B = 1.445647641019636;  C = 1.626987699542094;  D = 0.7210945179870265; E = 13.86384647080068;   G = 17.65757632934944 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
